# Auto-generated edit script: updates market-price derived columns (H-N)
# across multiple worksheets to match the target snapshot.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 992.5
$ws.Range("I19").Value = 588.3333
$ws.Range("J19").Value = 1235
$ws.Range("K19").Value = 588.3333
$ws.Range("L19").Value = 1235
$ws.Range("M19").Value = -413.3333
$ws.Range("N19").Value = -1585
$ws.Range("H86").Value = 58826140
$ws.Range("I86").Value = 100002770
$ws.Range("K86").Value = 100002770
$ws.Range("M86").Value = -100001647
$ws.Range("H89").Value = 58826140
$ws.Range("I89").Value = 100002770
$ws.Range("K89").Value = 500013850
$ws.Range("M89").Value = -500008234
$ws.Range("H99").Value = 333335330
$ws.Range("I99").Value = 3000
$ws.Range("K99").Value = 9000
$ws.Range("M99").Value = -7502
$ws.Range("H107").Value = 12346761
$ws.Range("I107").Value = 1008.05
$ws.Range("J107").Value = 47620340
$ws.Range("K107").Value = 1008.05
$ws.Range("L107").Value = 47620340
$ws.Range("M107").Value = 911.95
$ws.Range("N107").Value = -47624180
$ws.Range("H137").Value = 2088987.8
$ws.Range("I137").Value = 4867.5713
$ws.Range("K137").Value = 14602.7139
$ws.Range("M137").Value = -12052.7139
$ws.Range("H138").Value = 9165.134
$ws.Range("J138").Value = 3985.7144
$ws.Range("L138").Value = 11957.1432
$ws.Range("N138").Value = -22237.1432

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 859.5333000000001
$ws.Range("I2").Value = 822.53845
$ws.Range("K2").Value = 822.53845
$ws.Range("M2").Value = -709.53845
$ws.Range("H32").Value = 2581.6516
$ws.Range("I32").Value = 855.6712
$ws.Range("J32").Value = 10456.4375
$ws.Range("K32").Value = 855.6712
$ws.Range("L32").Value = 10456.4375
$ws.Range("M32").Value = -568.6712
$ws.Range("N32").Value = -11030.4375
$ws.Range("H61").Value = 2606496.8
$ws.Range("I61").Value = 88994.586
$ws.Range("K61").Value = 88994.586
$ws.Range("M61").Value = -88782.586
$ws.Range("H74").Value = 589215
$ws.Range("I74").Value = 2720.5908
$ws.Range("J74").Value = 1395644.9
$ws.Range("K74").Value = 2720.5908
$ws.Range("L74").Value = 1395644.9
$ws.Range("M74").Value = -1846.5908
$ws.Range("N74").Value = -1397392.9
$ws.Range("H77").Value = 589215
$ws.Range("I77").Value = 2720.5908
$ws.Range("J77").Value = 1395644.9
$ws.Range("K77").Value = 13602.954
$ws.Range("L77").Value = 6978224.5
$ws.Range("M77").Value = -9234.954
$ws.Range("N77").Value = -6986960.5
$ws.Range("H102").Value = 9355.154
$ws.Range("I102").Value = 11561.7
$ws.Range("K102").Value = 11561.7
$ws.Range("M102").Value = -9939.700000000001
$ws.Range("H116").Value = 859.5333000000001
$ws.Range("I116").Value = 822.53845
$ws.Range("K116").Value = 822.53845
$ws.Range("M116").Value = 1471.46155
$ws.Range("H123").Value = 142490
$ws.Range("J123").Value = 142490
$ws.Range("L123").Value = 142490
$ws.Range("N123").Value = -152290
$ws.Range("H132").Value = 1848.0555
$ws.Range("I132").Value = 1522.2307
$ws.Range("K132").Value = 4566.6921
$ws.Range("M132").Value = -2036.6921
$ws.Range("H135").Value = 63515.43
$ws.Range("J135").Value = 63515.43
$ws.Range("L135").Value = 63515.43
$ws.Range("N135").Value = -73655.42999999999
$ws.Range("H136").Value = 2606496.8
$ws.Range("I136").Value = 88994.586
$ws.Range("K136").Value = 266983.758
$ws.Range("M136").Value = -264433.758

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 859.5333000000001
$ws.Range("I3").Value = 822.53845
$ws.Range("K3").Value = 822.53845
$ws.Range("M3").Value = -708.53845
$ws.Range("H7").Value = 900
$ws.Range("I7").Value = 900
$ws.Range("K7").Value = 900
$ws.Range("M7").Value = -787
$ws.Range("H22").Value = 2937.5
$ws.Range("I22").Value = 2250
$ws.Range("K22").Value = 2250
$ws.Range("M22").Value = -2077
$ws.Range("H95").Value = 46327
$ws.Range("J95").Value = 46327
$ws.Range("L95").Value = 46327
$ws.Range("N95").Value = -51819
$ws.Range("H122").Value = 51934.195
$ws.Range("J122").Value = 51934.195
$ws.Range("L122").Value = 51934.195
$ws.Range("N122").Value = -61734.195
$ws.Range("H133").Value = 89990
$ws.Range("J133").Value = 89990
$ws.Range("L133").Value = 89990
$ws.Range("N133").Value = -100110
$ws.Range("H134").Value = 22502308
$ws.Range("I134").Value = 2114.0334
$ws.Range("K134").Value = 6342.100199999999
$ws.Range("M134").Value = -3807.100199999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1874
$ws.Range("I22").Value = 1423.8
$ws.Range("K22").Value = 1423.8
$ws.Range("M22").Value = -1073.8
$ws.Range("H97").Value = 52999.8
$ws.Range("J97").Value = 52999.8
$ws.Range("L97").Value = 52999.8
$ws.Range("N97").Value = -54981.8
$ws.Range("H105").Value = 1087.3478
$ws.Range("I105").Value = 1105.45
$ws.Range("K105").Value = 1105.45
$ws.Range("M105").Value = 641.55
$ws.Range("H107").Value = 1144.0303
$ws.Range("J107").Value = 1078.4166
$ws.Range("L107").Value = 1078.4166
$ws.Range("N107").Value = -4918.4166
$ws.Range("H123").Value = 169990
$ws.Range("J123").Value = 169990
$ws.Range("L123").Value = 169990
$ws.Range("N123").Value = -179790
$ws.Range("H132").Value = 37040796
$ws.Range("I132").Value = 3902.75
$ws.Range("J132").Value = 66670308
$ws.Range("K132").Value = 11708.25
$ws.Range("L132").Value = 200010924
$ws.Range("M132").Value = -9178.25
$ws.Range("N132").Value = -200015984
$ws.Range("H133").Value = 149999.67
$ws.Range("J133").Value = 149999.67
$ws.Range("L133").Value = 149999.67
$ws.Range("N133").Value = -155059.67
$ws.Range("H141").Value = 158139.08
$ws.Range("J141").Value = 186481.3
$ws.Range("L141").Value = 186481.3
$ws.Range("N141").Value = -196841.3

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H111").Value = 1105.4
$ws.Range("I111").Value = 1105.4
$ws.Range("K111").Value = 3316.2
$ws.Range("M111").Value = -249.2000000000003
$ws.Range("H113").Value = 729.5
$ws.Range("I113").Value = 196
$ws.Range("J113").Value = 836.2
$ws.Range("K113").Value = 588
$ws.Range("L113").Value = 2508.6
$ws.Range("M113").Value = 1582
$ws.Range("N113").Value = -6848.6
$ws.Range("H119").Value = 10830.667
$ws.Range("I119").Value = 4247.875
$ws.Range("K119").Value = 12743.625
$ws.Range("M119").Value = -7905.625
$ws.Range("H123").Value = 7400
$ws.Range("J123").Value = 10333.333
$ws.Range("L123").Value = 30999.999
$ws.Range("N123").Value = -35899.999
$ws.Range("H125").Value = 4995
$ws.Range("I125").Value = 5000
$ws.Range("K125").Value = 15000
$ws.Range("M125").Value = -10080
$ws.Range("H132").Value = 92195.63
$ws.Range("I132").Value = 1035.6
$ws.Range("J132").Value = 168162.33
$ws.Range("K132").Value = 9320.4
$ws.Range("L132").Value = 1513460.97
$ws.Range("M132").Value = -6790.4
$ws.Range("N132").Value = -1518520.97

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 1050000
$ws.Range("J7").Value = 1050000
$ws.Range("L7").Value = 1050000
$ws.Range("N7").Value = -1050224
$ws.Range("H8").Value = 1050000
$ws.Range("J8").Value = 1050000
$ws.Range("L8").Value = 1050000
$ws.Range("N8").Value = -1050278
$ws.Range("H46").Value = 32061.5
$ws.Range("J46").Value = 32061.5
$ws.Range("L46").Value = 32061.5
$ws.Range("N46").Value = -32373.5
$ws.Range("H95").Value = 35535
$ws.Range("J95").Value = 32380.334
$ws.Range("L95").Value = 32380.334
$ws.Range("N95").Value = -37872.334
$ws.Range("H97").Value = 804
$ws.Range("J97").Value = 871.3333
$ws.Range("L97").Value = 871.3333
$ws.Range("N97").Value = -1863.3333
$ws.Range("H113").Value = 3144.7222
$ws.Range("I113").Value = 2555.1538
$ws.Range("J113").Value = 4677.6
$ws.Range("K113").Value = 2555.1538
$ws.Range("L113").Value = 4677.6
$ws.Range("M113").Value = -385.1538
$ws.Range("N113").Value = -9017.6
$ws.Range("H126").Value = 2047.9
$ws.Range("I126").Value = 1420.5714
$ws.Range("K126").Value = 4261.7142
$ws.Range("M126").Value = -1791.7142

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 7640.7646
$ws.Range("I46").Value = 8653.429
$ws.Range("K46").Value = 8653.429
$ws.Range("M46").Value = -8465.429
$ws.Range("H122").Value = 3319.16
$ws.Range("J122").Value = 4595.8
$ws.Range("L122").Value = 13787.4
$ws.Range("N122").Value = -18687.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("M13").ClearContents()
$ws.Range("H19").Value = 451008
$ws.Range("J19").Value = 451008
$ws.Range("L19").Value = 451008
$ws.Range("N19").Value = -451356
$ws.Range("H136").Value = 4300.9443
$ws.Range("I136").Value = 4713.524
$ws.Range("K136").Value = 14140.572
$ws.Range("M136").Value = -11590.572
